$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaktiomaksut")

$ws.Range("A78").Value = "ei tiedossa"
$ws.Range("B78").Value = "C"

$ws.Range("B79").Select()
